# Weekly update: a new price record (week of 2023-03-24, date serial 45009)
# is inserted at the top of the Terminal Hortofrutícola Agro Chillán /
# Arándano (blue) block (row 33), pushing the existing rows 33-44 down to
# rows 34-45 and extending the used range to A1:T45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 33; this shifts rows 33:44 down to 34:45
# and Excel automatically grows the sheet dimension to A1:T45.
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new weekly record.
$ws.Range("A33").Value = 7
$ws.Range("B33").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C33").Value = "Ñuble"
$ws.Range("D33").Value = 45009
$ws.Range("E33").Value = 16
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100101
$ws.Range("H33").Value = "Berries"
$ws.Range("I33").Value = 100101001
$ws.Range("J33").Value = "Arándano (blue)"
$ws.Range("K33").Value = "Sin especificar"
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 30
$ws.Range("N33").Value = 4000
$ws.Range("O33").Value = 4000
$ws.Range("P33").Value = 4000
$ws.Range("Q33").Value = "`$/bandeja 2 kilos"
$ws.Range("R33").Value = "Provincia de Diguillín"
$ws.Range("S33").Value = 2000
$ws.Range("T33").Value = 2
